$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the selection on the "Program" sheet (was A13, becomes D1:F14)
# ---------------------------------------------------------------------------
$wsProgram = $wb.Worksheets.Item("Program")
[void]$wsProgram.Activate()
[void]$wsProgram.Range("D1:F14").Select()

# ---------------------------------------------------------------------------
# 2) Add a new "Msg" worksheet after the last existing sheet
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMsg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsMsg.Name = "Msg"

# ---------------------------------------------------------------------------
# 3) Populate the new sheet with the program-scenario / error-message table.
#    Cells are written in the specific order that reproduces the shared
#    string table ordering of the target workbook.
# ---------------------------------------------------------------------------
$wsMsg.Range("A2").Value = "Name Empty"
$wsMsg.Range("C2").Value = "Program name is required."

$wsMsg.Range("A3").Value = "Desc Empty"
$wsMsg.Range("C3").Value = "Description is required."

$wsMsg.Range("A4").Value = "Status Empty"
$wsMsg.Range("C4").Value = "Status is required."

$wsMsg.Range("A5").Value = "Name Spc char"
$wsMsg.Range("C5").Value = "This field should start with an alphabet, no special char and min 2 char."

$wsMsg.Range("A6").Value = "Name Strt char"
$wsMsg.Range("A7").Value = "Name Strt num"
$wsMsg.Range("A8").Value = "Name min char"
$wsMsg.Range("A9").Value = "Desc Strt num"
$wsMsg.Range("A10").Value = "Desc min char"

$wsMsg.Range("C1").Value = "errmsg"
$wsMsg.Range("A1").Value = "testscn"
$wsMsg.Range("B1").Value = "Input"

$wsMsg.Range("B5").Value = "Test#"
$wsMsg.Range("B6").Value = "#Test"
$wsMsg.Range("B7").Value = "01Test"
$wsMsg.Range("B8").Value = "T"

$wsMsg.Range("C9").Value = "This field should start with an alphabet and min 2 char."

$wsMsg.Range("C6").Value = "This field should start with an alphabet, no special char and min 2 char."
$wsMsg.Range("C7").Value = "This field should start with an alphabet, no special char and min 2 char."
$wsMsg.Range("C8").Value = "This field should start with an alphabet, no special char and min 2 char."
$wsMsg.Range("B9").Value = "01Test"
$wsMsg.Range("C10").Value = "This field should start with an alphabet and min 2 char."
$wsMsg.Range("B10").Value = "T"

# ---------------------------------------------------------------------------
# 4) Column widths for A and B
# ---------------------------------------------------------------------------
$wsMsg.Columns.Item(1).ColumnWidth = 16.5
$wsMsg.Columns.Item(2).ColumnWidth = 16.5

# ---------------------------------------------------------------------------
# 5) Final selection on the Msg sheet, which also becomes the active sheet/tab
# ---------------------------------------------------------------------------
[void]$wsMsg.Range("C10").Select()
